$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$products = @(
    "produit 101",
    "produit 102",
    "produit 103",
    "produit 104",
    "produit 105",
    "produit 106",
    "produit 107",
    "produit 108",
    "produit 109",
    "produit 110",
    "produit 111",
    "produit 112",
    "produit 113"
)

for ($i = 0; $i -lt $products.Length; $i++) {
    $row = 5 + $i
    $ws.Cells.Item($row, 1).Value = $products[$i]
}

$ws.Range("G15").Select()
